$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.980.76"
$ws.Range("E2").Value = "  -2.94%  "
$ws.Range("D3").Value = "3.369.11"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'568.33"
$ws.Range("D6").Value = "'148.90"
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'7.98"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("E11").Value = "  +2.22%  "
$ws.Range("D12").Value = "3.950.31"
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "'28.03"
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("D15").Value = "3.367.35"
$ws.Range("E15").Value = "  -2.38%  "
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "61.059.98"
$ws.Range("E17").Value = "  -2.94%  "
$ws.Range("D18").Value = "'6.36"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").Value = "'14.51"
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").Value = "'8.93"
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D21").Value = "'375.99"
$ws.Range("E21").Value = "  -3.27%  "
$ws.Range("D22").Value = "'75.42"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "'0.562"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "3.505.80"
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("D26").Value = "'0.0000108"
$ws.Range("E26").Value = "  -5.27%  "
$ws.Range("E27").Value = "  -3.86%  "
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  -3.70%  "
$ws.Range("D33").Value = "'22.91"
$ws.Range("E33").Value = "  -1.63%  "
$ws.Range("E34").Value = "  -3.31%  "
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("D36").Value = "'169.76"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("E37").Value = "  -3.43%  "
$ws.Range("E38").Value = "  -2.58%  "
$ws.Range("D39").Value = "'29.03"
$ws.Range("E39").Value = "  -9.31%  "
$ws.Range("D40").Value = "3.403.91"
$ws.Range("E40").Value = "  -2.35%  "
$ws.Range("D41").Value = "'0.0755"
$ws.Range("D42").Value = "'0.762"
$ws.Range("E42").Value = "  -3.79%  "
$ws.Range("E43").Value = "  -1.15%  "
$ws.Range("E44").Value = "  -2.55%  "
$ws.Range("E45").Value = "  -5.70%  "
$ws.Range("D46").Value = "2.492.68"
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("D47").Value = "'6.69"
$ws.Range("E47").Value = "  -2.89%  "
$ws.Range("D48").Value = "'22.65"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  -2.33%  "
$ws.Range("D51").Value = "'0.818"
$ws.Range("E51").Value = "  +0.35%  "
